$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update clinic address info
$ws.Range("AD2").Value = "phoenix"
$ws.Range("AC2").Value = "1201 s 7th ave"
$ws.Range("AE2").Value = 85007

# Update name fields
$ws.Range("W2").Value = "test"
$ws.Range("V2").Value = "test3"

# Update vaccine / assessment Yes/No answers
$ws.Range("H2").Value = "Yes"
$ws.Range("N2").Value = "No"
$ws.Range("P2").Value = "No"

# Reset view: scroll back to A1 and select P6
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P6").Select()
